$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Campaign Data" sub-table (rows 8-14) lost its rollup: the section
# header/category column and the Default/Automation breakdown columns are
# gone, leaving only the Total column (now all zero) and a blanked-out
# WP_DEMO column. ---
$ws.Range("A8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = ""

$ws.Range("A9").ClearContents()
$ws.Range("B9").Value = 0
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = ""

$ws.Range("A10").ClearContents()
$ws.Range("B10").Value = 0
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").Value = ""

$ws.Range("A11").ClearContents()
$ws.Range("B11").Value = 0
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = ""

$ws.Range("A12").ClearContents()
$ws.Range("B12").Value = 0
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").Value = ""

$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = 0
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = ""

$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = 0
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = ""

# --- Leads row: WP_DEMO workspace count dropped to 0 ---
$ws.Range("E17").Value = 0

# --- "Program Data" sub-table (rows 18-19): Total column cleared, WP_DEMO blanked ---
$ws.Range("B18").ClearContents()
$ws.Range("E18").Value = ""

$ws.Range("B19").ClearContents()
$ws.Range("E19").Value = ""

# --- Total WorkSpace count dropped to 0 ---
$ws.Range("B27").Value = 0

# --- New rows appended for the Design Studio / asset categories ---
$ws.Range("A30").Value = "Library"
$ws.Range("B30").Value = 3

$ws.Range("A31").Value = "Web Personalize"
$ws.Range("B31").Value = $false

$ws.Range("A32").Value = "Target Account Management"
$ws.Range("B32").Value = $true

$ws.Range("A33").Value = "Predictive Content"
$ws.Range("B33").Value = $false

# --- Selection moved to D31 (single cell) ---
$ws.Range("D31").Select()
